$wb = $excel.ActiveWorkbook

# --- "wheat" sheet: update row 2 and add a new row 3 ---
$wsWheat = $wb.Worksheets.Item("wheat")

$wsWheat.Range("A2").Value = "BSP"
$wsWheat.Range("B2").Value = "Gujarat"
$wsWheat.Range("C2").Value = "SMBX+FCSJ"
$wsWheat.Range("D2").Value = "Jammu & Kashmir"
$wsWheat.Range("E2").Value = "Wheat"
$wsWheat.Range("F2").Value = 1

$wsWheat.Range("A3").Value = "BH"
$wsWheat.Range("B3").Value = "Chattisgarh"
$wsWheat.Range("C3").Value = "BTI"
$wsWheat.Range("D3").Value = "Punjab"
$wsWheat.Range("E3").Value = "Wheat"
$wsWheat.Range("F3").Value = 1

# --- "rra" sheet: remove the data row, leaving only the header ---
$wsRra = $wb.Worksheets.Item("rra")
$wsRra.Range("A2:F2").Delete()
